$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (A1:H1) -------------------------------------------------
# Old headers: nome, placa, modelo, ano, cor, valor_compra, observacao, status, Data de Cadastro (A1:I1)
# New headers: nome, marca, cor, nome_dono, observacoes, valor_compra, status, Data de Cadastro (A1:H1)
$ws.Range("A1").Value = "nome"
$ws.Range("B1").Value = "marca"
$ws.Range("C1").Value = "cor"
$ws.Range("D1").Value = "nome_dono"
$ws.Range("E1").Value = "observacoes"
$ws.Range("F1").Value = "valor_compra"
$ws.Range("G1").Value = "status"
$ws.Range("H1").Value = "Data de Cadastro"

# ---- Data row (A2:H2) ----------------------------------------------------
$ws.Range("A2").Value = "Yamaha MT-07"
$ws.Range("B2").Value = "Yamaha"
$ws.Range("C2").Value = "preto"
$ws.Range("D2").Value = "Matheus"
$ws.Range("E2").Value = "teste"
$ws.Range("F2").Value = 15000
$ws.Range("G2").Value = "Ativo"
$ws.Range("H2").Value = "2024-05-17 16:22:23"

# ---- Drop the now-unused column I (was part of the old 9-column layout) --
$ws.Columns("I").Delete() | Out-Null

# ---- Re-apply the AutoFilter so its <autoFilter ref=.../> shrinks to H ---
# Toggling off then on with the new range rewrites the stored autofilter ref.
$ws.Range("A1:H2").AutoFilter() | Out-Null
$ws.Range("A1:H2").AutoFilter() | Out-Null

# ---- Fix up the _FilterDatabase defined name to match the new range ------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "ConfigMotos!_FilterDatabase") {
        $n.RefersTo = "='ConfigMotos'!`$A`$1:`$H`$2"
    }
}
